# Oklahoma City Thunder 2025-26 stats workbook update
# Adds "Rebounds" and "3PM" raw-data sheets (taking over the slots formerly
# occupied by "Avg Points"/"Avg Assists"), and re-creates "Avg Points" /
# "Avg Assists" at the end of the workbook, plus two brand new summary
# sheets "Avg Rebounds" and "Avg 3PM".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Shared column headers used by the four per-game raw-data sheets
# ---------------------------------------------------------------------
$headers = @(
    "Game Time (PST)", "Opponent", "Shai Gilgeous-Alexander", "Chris Youngblood",
    "Luguentz Dort", "Jaylin Williams", "Chet Holmgren", "Alex Caruso",
    "Isaiah Joe", "Ousmane Dieng", "Branden Carlson", "Aaron Wiggins",
    "Cason Wallace", "Brooks Barnhizer", "Ajay Mitchell", "Isaiah Hartenstein"
)

# Raw per-game Rebounds data: Date, Opponent, then 14 player columns (C..P)
$reboundsRows = @(
    @("2025-10-21","HOU",5,0,6,0,7,2,0,0,0,1,7,0,2,8),
    @("2025-10-23","IND",8,0,8,2,12,0,0,1,0,9,0,3,3,14),
    @("2025-10-25","ATL",5,1,2,3,11,0,0,2,0,5,3,1,7,10),
    @("2025-10-27","DAL",4,0,5,7,11,0,0,0,0,3,5,1,7,12),
    @("2025-10-28","SAC",9,0,5,5,0,2,0,0,1,6,1,1,0,14),
    @("2025-10-30","WAS",3,1,6,9,0,2,6,0,2,4,5,0,1,8),
    @("2025-11-02","NOP",2,1,0,3,0,2,5,1,5,1,0,3,5,14),
    @("2025-11-04","LAC",4,2,3,5,5,1,3,0,1,2,3,0,4,7),
    @("2025-11-05","POR",9,0,0,9,0,0,0,6,1,2,8,0,4,11),
    @("2025-11-07","SAC",1,0,0,3,7,2,6,5,0,0,2,2,5,19),
    @("2025-11-09","MEM",7,0,0,3,7,4,2,1,0,0,4,4,6,13)
)

# Raw per-game 3PM data: Date, Opponent, then 14 player columns (C..P)
$threePmRows = @(
    @("2025-10-21","HOU",1,0,0,0,2,2,0,0,0,2,3,0,3,0),
    @("2025-10-23","IND",2,1,1,1,0,0,0,0,0,5,0,0,0,0),
    @("2025-10-25","ATL",2,0,1,1,6,0,0,0,0,1,1,0,2,0),
    @("2025-10-27","DAL",1,0,2,0,2,0,0,0,0,1,0,0,0,0),
    @("2025-10-28","SAC",1,0,1,1,0,2,0,0,0,4,2,0,3,0),
    @("2025-10-30","WAS",3,0,1,2,0,0,5,0,1,0,2,0,2,0),
    @("2025-11-02","NOP",1,2,0,4,0,2,3,1,1,3,3,0,0,0),
    @("2025-11-04","LAC",4,0,1,1,0,2,6,0,0,2,2,0,0,0),
    @("2025-11-05","POR",3,0,0,1,0,0,3,0,0,7,1,0,0,0),
    @("2025-11-07","SAC",2,0,0,2,0,0,3,1,0,0,3,0,2,0),
    @("2025-11-09","MEM",3,0,0,0,1,1,2,0,0,0,1,0,2,0)
)

# Avg Points (Player, value) -- unchanged data, simply relocated to a new sheet
$avgPointsRows = @(
    @("Shai Gilgeous-Alexander", 33.18181818181818),
    @("Chet Holmgren", 19.14285714285714),
    @("Ajay Mitchell", 17.18181818181818),
    @("Isaiah Joe", 15.16666666666667),
    @("Aaron Wiggins", 14.77777777777778),
    @("Isaiah Hartenstein", 13.36363636363636),
    @("Cason Wallace", 9.199999999999999),
    @("Luguentz Dort", 7.142857142857143),
    @("Alex Caruso", 5.714285714285714),
    @("Jaylin Williams", 4.727272727272728),
    @("Branden Carlson", 2.571428571428572),
    @("Chris Youngblood", 1.5),
    @("Ousmane Dieng", 1.5),
    @("Brooks Barnhizer", 0.8888888888888888)
)

# Avg Assists (Player, value) -- unchanged data, simply relocated to a new sheet
$avgAssistsRows = @(
    @("Shai Gilgeous-Alexander", 6),
    @("Ajay Mitchell", 4.090909090909091),
    @("Isaiah Hartenstein", 3.272727272727273),
    @("Cason Wallace", 2.9),
    @("Alex Caruso", 2.571428571428572),
    @("Aaron Wiggins", 2.222222222222222),
    @("Jaylin Williams", 2),
    @("Chet Holmgren", 1.857142857142857),
    @("Luguentz Dort", 1.571428571428571),
    @("Isaiah Joe", 0.8333333333333334),
    @("Ousmane Dieng", 0.5),
    @("Branden Carlson", 0.2857142857142857),
    @("Brooks Barnhizer", 0.2222222222222222),
    @("Chris Youngblood", 0.1)
)

# Avg Rebounds (Player, value) -- new sheet
$avgReboundsRows = @(
    @("Isaiah Hartenstein", 11.81818181818182),
    @("Chet Holmgren", 8.571428571428571),
    @("Shai Gilgeous-Alexander", 5.181818181818182),
    @("Luguentz Dort", 5),
    @("Jaylin Williams", 4.454545454545454),
    @("Ajay Mitchell", 4),
    @("Cason Wallace", 3.8),
    @("Isaiah Joe", 3.666666666666667),
    @("Aaron Wiggins", 3.666666666666667),
    @("Alex Caruso", 2.142857142857143),
    @("Ousmane Dieng", 2),
    @("Brooks Barnhizer", 1.666666666666667),
    @("Branden Carlson", 1.428571428571429),
    @("Chris Youngblood", 0.5)
)

# Avg 3PM (Player, value) -- new sheet
$avg3pmRows = @(
    @("Isaiah Joe", 3.666666666666667),
    @("Aaron Wiggins", 2.777777777777778),
    @("Shai Gilgeous-Alexander", 2.090909090909091),
    @("Cason Wallace", 1.8),
    @("Chet Holmgren", 1.571428571428571),
    @("Alex Caruso", 1.285714285714286),
    @("Ajay Mitchell", 1.272727272727273),
    @("Jaylin Williams", 1.181818181818182),
    @("Luguentz Dort", 1),
    @("Chris Youngblood", 0.3),
    @("Branden Carlson", 0.2857142857142857),
    @("Ousmane Dieng", 0.25),
    @("Brooks Barnhizer", 0),
    @("Isaiah Hartenstein", 0)
)

# ---------------------------------------------------------------------
# Helper: write a "Game Time (PST)/Opponent/<14 players>" raw data table
# into a worksheet, replacing whatever it previously contained. Column A
# holds YYYY-MM-DD-looking text that Excel would otherwise auto-convert
# into a date serial, so it is force-formatted as Text first and then
# reset back to the default "Normal" style (matching the source data,
# which stores these as plain strings with no special number format).
# ---------------------------------------------------------------------
function Write-RawTable($sheet, $rows) {
    $sheet.Cells.Clear()
    for ($c = 0; $c -lt $headers.Length; $c++) {
        $sheet.Cells.Item(1, $c + 1).Value = $headers[$c]
    }
    for ($r = 0; $r -lt $rows.Length; $r++) {
        $row = $rows[$r]
        $dateCell = $sheet.Cells.Item($r + 2, 1)
        $dateCell.NumberFormat = "@"
        $dateCell.Value = $row[0]
        $dateCell.Style = "Normal"
        for ($c = 1; $c -lt $row.Length; $c++) {
            $sheet.Cells.Item($r + 2, $c + 1).Value = $row[$c]
        }
    }
}

# ---------------------------------------------------------------------
# Helper: write a "Player/<metric>" two-column summary table into a sheet.
# ---------------------------------------------------------------------
function Write-SummaryTable($sheet, $metricName, $rows) {
    $sheet.Cells.Clear()
    $sheet.Cells.Item(1, 1).Value = "Player"
    $sheet.Cells.Item(1, 2).Value = $metricName
    for ($r = 0; $r -lt $rows.Length; $r++) {
        $row = $rows[$r]
        $sheet.Cells.Item($r + 2, 1).Value = $row[0]
        $sheet.Cells.Item($r + 2, 2).Value = $row[1]
    }
}

# Reference header cell (bold / centered / bordered) used to carry over the
# same formatting onto every header row we (re)write below.
$headerStyleSource = $wb.Worksheets.Item(1).Range("A1")

# ---------------------------------------------------------------------
# 1) Turn the current sheet 3 ("Avg Points") into "Rebounds" raw data.
# 2) Turn the current sheet 4 ("Avg Assists") into "3PM" raw data.
# (sheetId / r:id values 3 and 4 are preserved this way.)
# ---------------------------------------------------------------------
$reboundsSheet = $wb.Worksheets.Item(3)
Write-RawTable $reboundsSheet $reboundsRows
$reboundsSheet.Name = "Rebounds"
$headerStyleSource.Copy()
$reboundsSheet.Range("A1:P1").PasteSpecial(-4122)

$threePmSheet = $wb.Worksheets.Item(4)
Write-RawTable $threePmSheet $threePmRows
$threePmSheet.Name = "3PM"
$headerStyleSource.Copy()
$threePmSheet.Range("A1:P1").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Append the four summary sheets, in order, at the end of the workbook:
# Avg Points, Avg Assists, Avg Rebounds, Avg 3PM.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$avgPointsSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$avgPointsSheet.Name = "Avg Points"
Write-SummaryTable $avgPointsSheet "Avg Points" $avgPointsRows
$headerStyleSource.Copy()
$avgPointsSheet.Range("A1:B1").PasteSpecial(-4122)

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$avgAssistsSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$avgAssistsSheet.Name = "Avg Assists"
Write-SummaryTable $avgAssistsSheet "Avg Assists" $avgAssistsRows
$headerStyleSource.Copy()
$avgAssistsSheet.Range("A1:B1").PasteSpecial(-4122)

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$avgReboundsSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$avgReboundsSheet.Name = "Avg Rebounds"
Write-SummaryTable $avgReboundsSheet "Avg Rebounds" $avgReboundsRows
$headerStyleSource.Copy()
$avgReboundsSheet.Range("A1:B1").PasteSpecial(-4122)

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$avg3pmSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$avg3pmSheet.Name = "Avg 3PM"
Write-SummaryTable $avg3pmSheet "Avg 3PM" $avg3pmRows
$headerStyleSource.Copy()
$avg3pmSheet.Range("A1:B1").PasteSpecial(-4122)

# Leave the first sheet selected/active, matching the original workbook.
$wb.Worksheets.Item(1).Select()
